# Add the new firmware-version rows (V.1.040 / 102 and V.1.050 / 103)
# to the "version" worksheet, mirroring the formatting of the most
# recent existing row (row 8) for every column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("version")

$xlPasteFormats = -4122

function Copy-RowFormat($srcRow, $dstRow) {
    foreach ($col in "A", "B", "C", "D", "E", "F") {
        $ws.Range("$col$srcRow").Copy() | Out-Null
        $ws.Range("$col$dstRow").PasteSpecial($xlPasteFormats) | Out-Null
    }
}

Copy-RowFormat 8 10
Copy-RowFormat 8 11

# Row 10 - V.1.040 / FW 102
$ws.Range("F10").Value = "fixed the problem of the spurious 1 second offline Alarm"
$ws.Range("A10").Value = "'102"
$ws.Range("D10").Value = "V.1.040"
$ws.Range("B10").Value = "'101"
$ws.Range("C10").Value = 971
$ws.Range("E10").Value = 44379

# Row 11 - V.1.050 / FW 103
$ws.Range("A11").Value = "'103"
$ws.Range("D11").Value = "V.1.050"
$ws.Range("F11").Value = "release from PLCM 9718"
$ws.Range("B11").Value = "'101"
$ws.Range("C11").Value = 980
$ws.Range("E11").Value = 44463

# Move the active selection like the authored workbook (F11 selected).
$ws.Range("F11").Select()
